$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update hardcoded AE/AF (date) input cells on rows 8-20.
# These feed the WORKDAY()-based formula chains elsewhere in the row, so
# every dependent formula cell (F,G,K,L,P,Q,Z,AA,AJ,AK,AN, etc.) recalculates
# automatically once these inputs change.
$ws.Range("AE8").Value = "12/09/2024"
$ws.Range("AF8").Value = "01/03/2025"

$ws.Range("AE9").Value = "06/04/2025"
$ws.Range("AF9").Value = "07/01/2025"

$ws.Range("AE10").Value = "07/07/2025"
$ws.Range("AF10").Value = "08/01/2025"

$ws.Range("AE11").Value = "08/05/2025"
$ws.Range("AF11").Value = "09/01/2025"

$ws.Range("AE12").Value = "08/05/2025"
$ws.Range("AF12").Value = "09/01/2025"

$ws.Range("AE13").Value = "08/05/2025"
$ws.Range("AF13").Value = "09/01/2025"

$ws.Range("AE14").Value = "10/07/2025"
$ws.Range("AF14").Value = "11/03/2025"

$ws.Range("AE15").Value = "10/07/2025"
$ws.Range("AF15").Value = "11/03/2025"

$ws.Range("AE16").Value = "11/04/2025"

$ws.Range("AE17").Value = "12/05/2025"
$ws.Range("AF17").Value = "01/01/2026"

$ws.Range("AE18").Value = "12/05/2025"
$ws.Range("AF18").Value = "01/01/2026"

$ws.Range("AE19").Value = "01/06/2026"
$ws.Range("AF19").Value = "02/01/2026"

$ws.Range("AE20").Value = "02/03/2026"
$ws.Range("AF20").Value = "03/01/2026"
